$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.356.65'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.269.76'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.52%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.38'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.632'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.49'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +5.62%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.642'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +4.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.47'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '59.17'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0963'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.29'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.42%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.603.44'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.81'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.877'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.267.25'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.342.75'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0987'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.27'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.76'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.27'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +8.91%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '234.55'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.92'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.50'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.43'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.64'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '166.43'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.94'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.42'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +10.83%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.126'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.80%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '31.50'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +20.98%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0797'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.90%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.77'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +11.00%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.125'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.73'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0307'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.54'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +13.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.32'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.97'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +5.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.210'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +7.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.24'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.49'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.89'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -6.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.103'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.41%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.17'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.18'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.35%  '
